# This script applies the "Updated cryptos list" data refresh to Sheet1.
# It updates the Price (column D) and Volume(1h) (column E) values for
# most rows, and also reflects two rows whose coin (and all associated
# data) moved position in the source ranking:
#   - Maker/Hedera swapped positions in rows 34/35
#   - FirstDigitalUSD/Stellar swapped positions in rows 48/49

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for most rows.
# D-column values are forced to Text format first so Excel stores them
# as the exact original text (preserving trailing zeros / dot-grouping)
# instead of auto-converting them to floating point numbers.
$ws.Range("D2").NumberFormat = '@'
$ws.Range("D2").Value = '67.814.87'
$ws.Range("E2").Value = '  +0.26%  '
$ws.Range("D3").NumberFormat = '@'
$ws.Range("D3").Value = '3.304.07'
$ws.Range("E3").Value = '  -2.10%  '
$ws.Range("D4").NumberFormat = '@'
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").NumberFormat = '@'
$ws.Range("D5").Value = '581.99'
$ws.Range("E5").Value = '  -1.57%  '
$ws.Range("D6").NumberFormat = '@'
$ws.Range("D6").Value = '175.67'
$ws.Range("E6").Value = '  -6.50%  '
$ws.Range("D7").NumberFormat = '@'
$ws.Range("D7").Value = '0.998'
$ws.Range("E7").Value = '  -0.13%  '
$ws.Range("D8").NumberFormat = '@'
$ws.Range("D8").Value = '0.580'
$ws.Range("E8").Value = '  -3.00%  '
$ws.Range("D9").NumberFormat = '@'
$ws.Range("D9").Value = '3.299.53'
$ws.Range("E9").Value = '  -1.93%  '
$ws.Range("D10").NumberFormat = '@'
$ws.Range("D10").Value = '0.174'
$ws.Range("E10").Value = '  -5.44%  '
$ws.Range("D11").NumberFormat = '@'
$ws.Range("D11").Value = '0.572'
$ws.Range("E11").Value = '  -2.69%  '
$ws.Range("D12").NumberFormat = '@'
$ws.Range("D12").Value = '45.39'
$ws.Range("E12").Value = '  -4.55%  '
$ws.Range("D13").NumberFormat = '@'
$ws.Range("D13").Value = '0.0000268'
$ws.Range("E13").Value = '  -2.57%  '
$ws.Range("D14").NumberFormat = '@'
$ws.Range("D14").Value = '665.94'
$ws.Range("E14").Value = '  +3.80%  '
$ws.Range("D15").NumberFormat = '@'
$ws.Range("D15").Value = '3.821.58'
$ws.Range("E15").Value = '  -2.39%  '
$ws.Range("D16").NumberFormat = '@'
$ws.Range("D16").Value = '8.33'
$ws.Range("E16").Value = '  -3.52%  '
$ws.Range("D17").NumberFormat = '@'
$ws.Range("D17").Value = '67.716.42'
$ws.Range("E17").Value = '  +0.19%  '
$ws.Range("D19").NumberFormat = '@'
$ws.Range("D19").Value = '3.287.21'
$ws.Range("E19").Value = '  -2.51%  '
$ws.Range("D20").NumberFormat = '@'
$ws.Range("D20").Value = '17.38'
$ws.Range("E20").Value = '  -3.72%  '
$ws.Range("D21").NumberFormat = '@'
$ws.Range("D21").Value = '10.83'
$ws.Range("E21").Value = '  -3.25%  '
$ws.Range("D22").NumberFormat = '@'
$ws.Range("D22").Value = '0.884'
$ws.Range("E22").Value = '  -2.98%  '
$ws.Range("D23").NumberFormat = '@'
$ws.Range("D23").Value = '5.44'
$ws.Range("E23").Value = '  +6.24%  '
$ws.Range("D24").NumberFormat = '@'
$ws.Range("D24").Value = '17.08'
$ws.Range("E24").Value = '  -5.37%  '
$ws.Range("D25").NumberFormat = '@'
$ws.Range("D25").Value = '98.31'
$ws.Range("E25").Value = '  -1.50%  '
$ws.Range("D26").NumberFormat = '@'
$ws.Range("D26").Value = '3.85'
$ws.Range("E26").Value = '  -4.06%  '
$ws.Range("D27").NumberFormat = '@'
$ws.Range("D27").Value = '2.66'
$ws.Range("E27").Value = '  -6.98%  '
$ws.Range("D28").NumberFormat = '@'
$ws.Range("D28").Value = '9.19'
$ws.Range("E28").Value = '  -5.73%  '
$ws.Range("D29").NumberFormat = '@'
$ws.Range("D29").Value = '32.78'
$ws.Range("E29").Value = '  +0.64%  '
$ws.Range("D30").NumberFormat = '@'
$ws.Range("D30").Value = '8.36'
$ws.Range("E30").Value = '  -4.17%  '
$ws.Range("D31").NumberFormat = '@'
$ws.Range("D31").Value = '7.03'
$ws.Range("E31").Value = '  +1.16%  '
$ws.Range("D32").NumberFormat = '@'
$ws.Range("D32").Value = '581.19'
$ws.Range("E32").Value = '  -5.30%  '
$ws.Range("D33").NumberFormat = '@'
$ws.Range("D33").Value = '10.93'
$ws.Range("E33").Value = '  -1.89%  '
$ws.Range("D36").NumberFormat = '@'
$ws.Range("D36").Value = '1.00'
$ws.Range("E36").Value = '  +0.11%  '
$ws.Range("D37").NumberFormat = '@'
$ws.Range("D37").Value = '3.37'
$ws.Range("E37").Value = '  -13.00%  '
$ws.Range("D38").NumberFormat = '@'
$ws.Range("D38").Value = '55.58'
$ws.Range("E38").Value = '  -0.79%  '
$ws.Range("D39").NumberFormat = '@'
$ws.Range("D39").Value = '0.131'
$ws.Range("E39").Value = '  -1.24%  '
$ws.Range("D40").NumberFormat = '@'
$ws.Range("D40").Value = '32.39'
$ws.Range("E40").Value = '  -4.31%  '
$ws.Range("D41").NumberFormat = '@'
$ws.Range("D41").Value = '2.63'
$ws.Range("E41").Value = '  -8.15%  '
$ws.Range("D42").NumberFormat = '@'
$ws.Range("D42").Value = '3.04'
$ws.Range("E42").Value = '  -7.60%  '
$ws.Range("D43").NumberFormat = '@'
$ws.Range("D43").Value = '0.0₃0662'
$ws.Range("E43").Value = '  -6.57%  '
$ws.Range("D44").NumberFormat = '@'
$ws.Range("D44").Value = '3.25'
$ws.Range("E44").Value = '  -4.95%  '
$ws.Range("D45").NumberFormat = '@'
$ws.Range("D45").Value = '0.328'
$ws.Range("E45").Value = '  -5.16%  '
$ws.Range("D46").NumberFormat = '@'
$ws.Range("D46").Value = '0.0401'
$ws.Range("E46").Value = '  -5.23%  '
$ws.Range("D51").NumberFormat = '@'
$ws.Range("D51").Value = '2.75'
$ws.Range("E51").Value = '  -2.08%  '

# Update Volume(1h)-only rows
$ws.Range("E18").Value = '  -0.58%  '
$ws.Range("E47").Value = '  -0.55%  '
$ws.Range("E50").Value = '  -1.38%  '

# Row 34/35: Maker and Hedera swap places with updated values
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = '@'
$ws.Range("D34").Value = "0.103"
$ws.Range("E34").Value = '  -3.17%  '
$ws.Range("B35").Value = "Maker"
$ws.Range("C35").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D35").NumberFormat = '@'
$ws.Range("D35").Value = "3.750.73"
$ws.Range("E35").Value = '  -4.54%  '

# Row 48/49: FirstDigitalUSD and Stellar swap places with updated values
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").NumberFormat = '@'
$ws.Range("D48").Value = "0.127"
$ws.Range("E48").Value = '  -2.47%  '
$ws.Range("B49").Value = "FirstDigitalUSD"
$ws.Range("C49").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D49").NumberFormat = '@'
$ws.Range("D49").Value = "1.00"
$ws.Range("E49").Value = '  +0.03%  '
